$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename existing "High Priority break-up" sheet and add the new sheet
#    "Major update - High Priority " right after it.
# ---------------------------------------------------------------------------
$wsOldBreakup = $wb.Worksheets.Item("High Priority break-up")
$wsOldBreakup.Name = "Interannual update - High Pri"

$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsOldBreakup)
$wsNew.Name = "Major update - High Priority "

# ---------------------------------------------------------------------------
# 2. Trends Status sheet updates
# ---------------------------------------------------------------------------
$wsTrends = $wb.Worksheets.Item("Trends Status")

# Row 2 - Rapid Decline
$wsTrends.Range("C2").Value = 0
$wsTrends.Range("D2:E2").ClearContents()

# Row 3 - Decline
$wsTrends.Range("B3").Value = 0
$wsTrends.Range("D3:E3").ClearContents()

# Row 4 - Stable
$wsTrends.Range("C4").Value = 0
$wsTrends.Range("D4:E4").ClearContents()

# Row 5 - Increase
$wsTrends.Range("D5:E5").ClearContents()

# Row 6 - Rapid Increase
$wsTrends.Range("B6").Value = 0
$wsTrends.Range("D6:E6").ClearContents()

# Row 7 - Trend Inconclusive
$wsTrends.Range("B7").Value = 9
$wsTrends.Range("C7").Value = 29

# Row 8 - Insufficient Data
$wsTrends.Range("C8").Value = 424

# ---------------------------------------------------------------------------
# 3. Priority Status sheet updates
# ---------------------------------------------------------------------------
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# 4. Species qualification sheet updates
# ---------------------------------------------------------------------------
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("A2").Value = "SoIB Assessment"
$wsSpecies.Range("B2").Value = 453
$wsSpecies.Range("B3").Value = 9
$wsSpecies.Range("C3").Value = 0
$wsSpecies.Range("B4").Value = 29
$wsSpecies.Range("C4").Value = 0

# ---------------------------------------------------------------------------
# 5. "Interannual update - High Pri" sheet (formerly "High Priority break-up")
#    gets new values for its two data rows.
# ---------------------------------------------------------------------------
$wsInterannual = $wb.Worksheets.Item("Interannual update - High Pri")
$wsInterannual.Range("B2").Value = 87
$wsInterannual.Range("C2").Value = 84.5
$wsInterannual.Range("D2").Value = 87
$wsInterannual.Range("E2").Value = 88.8

$wsInterannual.Range("B3").Value = 16
$wsInterannual.Range("C3").Value = 15.5
$wsInterannual.Range("D3").Value = 11
$wsInterannual.Range("E3").Value = 11.2

# ---------------------------------------------------------------------------
# 6. New "Major update - High Priority " sheet gets the content that used
#    to live in the original "High Priority break-up" sheet.
# ---------------------------------------------------------------------------
$wsNew.Range("A1").Value = "Break-up"
$wsNew.Range("B1").Value = "High Species (no.)"
$wsNew.Range("C1").Value = "High Species (perc.)"
$wsNew.Range("D1").Value = "New High Species (no.)"
$wsNew.Range("E1").Value = "New High Species (perc.)"

$wsNew.Range("A2").Value = "Trend New"
$wsNew.Range("B2").Value = 2
$wsNew.Range("C2").Value = 22.2
$wsNew.Range("D2").Value = 2
$wsNew.Range("E2").Value = 22.2

$wsNew.Range("A3").Value = "IUCN"
$wsNew.Range("B3").Value = 7
$wsNew.Range("C3").Value = 77.8
$wsNew.Range("D3").Value = 7
$wsNew.Range("E3").Value = 77.8

# Match the bold, centered header styling used by the other summary sheets.
$wsNew.Range("A1:E1").Font.Bold = $true
$wsNew.Range("A1:E1").HorizontalAlignment = -4108

$wsTrends.Select()
